# Rebuild the "total" sheet with the new expense-tracker dataset.
#
# Strategy: insert a brand-new worksheet, populate it completely with the
# new header row + 8 data rows (with per-row, non-shared formulas), then
# delete the old "total" sheet and rename the new one to "total". This
# naturally reproduces the sheetId bump, the loss of the stale
# _xlnm._FilterDatabase defined name (it lived on the old sheet), and the
# reset sheet view (no custom zoom) that show up in the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add()

# ---- headers (row 1) ----
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "expense_category"
$ws.Range("C1").Value = "expense_type"
$ws.Range("D1").Value = "value"
$ws.Range("E1").Value = "month"
$ws.Range("F1").Value = "year"
$ws.Range("G1").Value = "weekday_number"
$ws.Range("H1").Value = "weekday_text"
$ws.Range("I1").Value = "months_text"
$ws.Range("J1").Value = "store"
$ws.Range("K1").Value = "city"

# ---- data rows (2..9) ----
# row -> category, expense_type, value (number or formula string)
$data = @(
    @{ Row = 2; Category = "food";               Item = "zucchini";       Value = 5.3 },
    @{ Row = 3; Category = "food";               Item = "watermelon";     Value = 0.79 },
    @{ Row = 4; Category = "cleaning products";   Item = "window cleaner"; Value = 0.95 },
    @{ Row = 5; Category = "cleaning products";   Item = "detergent";      Value = "=1.19/2" },
    @{ Row = 6; Category = "food";               Item = "butter";         Value = 0.99 },
    @{ Row = 7; Category = "food";               Item = "cocoa";          Value = 1.05 },
    @{ Row = 8; Category = "transportation";      Item = "train ticket";   Value = 25 },
    @{ Row = 9; Category = "transportation";      Item = "bus ticket";     Value = 45 }
)

foreach ($d in $data) {
    $r = $d.Row

    $ws.Cells.Item($r, 1).Value = 43831
    $ws.Cells.Item($r, 2).Value = $d.Category
    $ws.Cells.Item($r, 3).Value = $d.Item

    if ($d.Value -is [string]) {
        $ws.Cells.Item($r, 4).Formula = $d.Value
    } else {
        $ws.Cells.Item($r, 4).Value = $d.Value
    }

    $ws.Cells.Item($r, 5).Formula = "=MONTH(A$r)"
    $ws.Cells.Item($r, 6).Formula = "=YEAR(A$r)"
    $ws.Cells.Item($r, 7).Formula = "=WEEKDAY(A$r, 2)"
    $ws.Cells.Item($r, 8).Formula = '=CHOOSE(WEEKDAY(A' + $r + ', 2), "Monday", "Tuesday","Wednesday", "Thursday", "Friday", "Saturday","Sunday")'
    $ws.Cells.Item($r, 9).Formula = '=TEXT(A' + $r + ', "MMM")'
}

# Date formatting for A2:A9 - "m/d/yy" maps to the built-in short-date
# numFmtId (14), matching the style used by the original sheet.
$ws.Range("A2:A9").NumberFormat = "m/d/yy"

# Column widths: match Excel's own best-fit behaviour for the populated
# data (mirrors the bestFit/customWidth columns in the target sheet).
$ws.Columns.AutoFit()

# ---- swap sheets: drop the old "total" and promote the new one ----
# (fetched fresh by name - handles are index-anchored, and inserting the
# new sheet at position 1 shifted the original "total" to position 2)
$wb.Sheets.Item("total").Delete()
$ws.Name = "total"

# ---- view state: select H10 on the (now single) "total" sheet ----
$ws.Range("H10").Select()
